$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Replace the old "lab-dashboard-testOffered" entry (row 35) with a new
# "/api/product" API block that spans rows 35:39, matching the layout of the
# other API-description blocks above it (merged API-name cell in column B,
# endpoint/method/payload/return/protected columns to the right).
# ---------------------------------------------------------------------------

# --- Row 35: header row of the new block -----------------------------------
$ws.Range("B35").Value = "/api/product"

$productJson = @"
{
  "labRef": "labID"
  "productList": [
    {
      "name": "Blood Test",
      "description": "Complete blood count",
      ...
    },
    {
      "name": "Urine Analysis",
      "description": "Urinalysis test",
    ...
    },
    {
      "name": "X-Ray",
      "description": "Chest X-ray",
     ......
    }
  ]
}
"@
$ws.Range("E35").Value = $productJson

$ws.Range("F35").Value = "updatedProduct JSON with ID"
$ws.Range("G35").Value = "Yes, opened as of now"

$ws.Rows(35).RowHeight = 300

# --- Row 36: delete ----------------------------------------------------
$ws.Range("C36").Value = "/delete/:id"
$ws.Range("D36").Value = "delete"
$ws.Range("G36").Value = "YES"

# --- Row 37: update (post) ----------------------------------------------
$ws.Range("C37").Value = "/update/:id"
$ws.Range("D37").Value = "post"
$ws.Range("E37").Value = "same as create"
$ws.Range("G37").Value = "YES"

# --- Row 38: get by id ----------------------------------------------------
$ws.Range("C38").Value = "/get/:id"
$ws.Range("D38").Value = "get"
$ws.Range("G38").Value = "YES"

# --- Row 39: trailing blank row of the block -------------------------------
# (no content, just part of the merged/bordered block)

# ---------------------------------------------------------------------------
# Formatting: apply the same look used by the existing blocks (B4:B8 etc.)
# to the new rows, then merge + highlight the API-name column.
# ---------------------------------------------------------------------------

# Whole detail block (C36:G39) gets a thin border around every cell.
$ws.Range("C36:G39").Borders.LineStyle = 1

# Endpoint / method / protected columns are centered (rows 36-38 only).
$ws.Range("C36:D38").HorizontalAlignment = -4108
$ws.Range("C36:D38").VerticalAlignment = -4108
$ws.Range("G36:G38").HorizontalAlignment = -4108
$ws.Range("G36:G38").VerticalAlignment = -4108

# Row 39 (the blank filler row) stays left/top aligned like row 3 - no
# explicit alignment override needed (default).

# Merge the API name cell across the whole block and style it like the
# other merged "API name" cells, then highlight it green.
$ws.Range("B35:B39").Merge()
$ws.Range("B35:B39").Borders.LineStyle = 1
$ws.Range("B35:B39").HorizontalAlignment = -4108
$ws.Range("B35:B39").VerticalAlignment = -4108
$ws.Range("B35:B39").Interior.Color = 5296274

# G35 (protected=Yes) also gets the green highlight, wrapped text, centered.
$ws.Range("G35").Interior.Color = 5296274
$ws.Range("G35").HorizontalAlignment = -4108
$ws.Range("G35").VerticalAlignment = -4108
$ws.Range("G35").WrapText = $true

# F35 (return) uses the same plain bordered/wrapped style as E10/E28's base
# font (non-bold "family 2" Calibri variant), matching the other payload
# cells that contain mixed-formatting rich text.
$ws.Range("F35").WrapText = $true
$ws.Range("F35").Borders.LineStyle = 1

# E35 (payload) keeps the plain wrap-text style used by the old cell.
$ws.Range("E35").WrapText = $true

# ---------------------------------------------------------------------------
# Reflect where the user ended up after the edit (scrolled down to the new
# block, with J35 selected).
# ---------------------------------------------------------------------------
$ws.Range("J35").Select()
$excel.ActiveWindow.ScrollRow = 35
